$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 167, shifting rows 167:216 down to 168:217.
$ws.Rows.Item(167).EntireRow.Insert()

# Populate the newly inserted row 167 with the new record
# (same Mercado/Región/Categoría/Unidad/Origen/Kg as the rest of the block).
$ws.Range("A167").Value = 2
$ws.Range("B167").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C167").Value = "Coquimbo"
$ws.Range("D167").Value = 44588
$ws.Range("E167").Value = 4
$ws.Range("F167").Value = 100112021
$ws.Range("G167").Value = "Ají"
$ws.Range("H167").Value = "Americana (o)"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 300
$ws.Range("K167").Value = 8000
$ws.Range("L167").Value = 10000
$ws.Range("M167").Value = 9000
$ws.Range("N167").Value = "`$/caja 25 kilos"
$ws.Range("O167").Value = "Provincia de Limarí"
$ws.Range("P167").Value = 360
$ws.Range("Q167").Value = 25
$ws.Range("R167").Value = "Hortaliza"
